# Generate Report for Archive
#
# Updates localization status from "Ready for handoff" to "In Translation"
# across the Overview summary sheet and each per-locale detail sheet, and
# tightens the affected "Status" column widths to match the shorter text
# (mirrors Excel auto-fitting those columns after the text change).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRange = $wsOverview.Range("E2:F3")
for ($r = 1; $r -le $overviewRange.Rows.Count; $r++) {
    for ($c = 1; $c -le $overviewRange.Columns.Count; $c++) {
        $cell = $overviewRange.Cells.Item($r, $c)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale sheets: Status column (C) ---
$localeSheets = @("zh-cn", "de-de")
foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C3")
    for ($r = 1; $r -le $statusRange.Rows.Count; $r++) {
        $cell = $statusRange.Cells.Item($r, 1)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
